# Correct omissions on license/acknowledgements slide
#
# Target slide: "License, Citation and Acknowledgements" -- the 2nd slide
# in the deck (ppt/slides/slide2.xml).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)          # "Content Placeholder 2"
$tr = $sh.TextFrame.TextRange

$curly = [char]0x2019            # RIGHT SINGLE QUOTATION MARK used in the deck

# ---------------------------------------------------------------------
# 1) Citation paragraph: add "August 2020, " before "online. DOI: "
# ---------------------------------------------------------------------
$full = $tr.Text
$searchOld = ", Software Productivity Track, in Argonne Training Program for Extreme Scale Computing (ATPESC), online. DOI: "
$writeNew  = ", Software Productivity Track, in Argonne Training Program for Extreme Scale Computing (ATPESC), August 2020, online. DOI: "
$idx = $full.IndexOf($searchOld)
if ($idx -ge 0) {
    $chars = $tr.Characters($idx + 1, $searchOld.Length)
    $chars.Text = $writeNew
}

# ---------------------------------------------------------------------
# 2) "Acknowledgements" heading paragraph: add 8pt space-before, keeping
#    <a:buNone/> after <a:spcBef> (toggle the bullet type to force the
#    paragraph properties to be re-emitted in that order).
# ---------------------------------------------------------------------
$tr = $sh.TextFrame.TextRange
$paraAck = $tr.Paragraphs(5, 1)
$pf = $paraAck.ParagraphFormat
$pf.SpaceBefore = 8
$pf.Bullet.Type = 1
$pf.Bullet.Type = 0

# ---------------------------------------------------------------------
# 3) Contributors paragraph: add ", Deborah Stevens" after "David Rogers"
# ---------------------------------------------------------------------
$tr = $sh.TextFrame.TextRange
$full = $tr.Text
$searchOld = ", Jared O'Neal, David Rogers"
$writeNew  = ", Jared O" + $curly + "Neal, David Rogers, Deborah Stevens"
$idx = $full.IndexOf($searchOld)
if ($idx -ge 0) {
    $chars = $tr.Characters($idx + 1, $searchOld.Length)
    $chars.Text = $writeNew
}

# ---------------------------------------------------------------------
# 4) Sandia acknowledgement paragraph: drop the stray "SAND NO " prefix
# ---------------------------------------------------------------------
$tr = $sh.TextFrame.TextRange
$full = $tr.Text
$searchOld = "This work was performed in part at Sandia National Laboratories. Sandia National Laboratories is a multi-mission laboratory managed and operated by National Technology and Engineering Solutions of Sandia, LLC., a wholly owned subsidiary of Honeywell International, Inc., for the U.S. Department of Energy's National Nuclear Security Administration under contract DE-NA0003525. SAND NO SAND2020-7957 PE"
$writeNew  = "This work was performed in part at Sandia National Laboratories. Sandia National Laboratories is a multi-mission laboratory managed and operated by National Technology and Engineering Solutions of Sandia, LLC., a wholly owned subsidiary of Honeywell International, Inc., for the U.S. Department of Energy" + $curly + "s National Nuclear Security Administration under contract DE-NA0003525. SAND2020-7957 PE"
$idx = $full.IndexOf($searchOld)
if ($idx -ge 0) {
    $chars = $tr.Characters($idx + 1, $searchOld.Length)
    $chars.Text = $writeNew
}
